{"js": "// Minor change: insert \"2  \" (a grammar-flagged interjection) right\n// before \"text\" in the second paragraph, turning\n//   \"This is test text with test idea.\"\n// into\n//   \"This is test 2  text with test idea.\"\n// and splitting the run the same way Word's grammar checker would\n// (wrapping the inserted word + the following word in a gramStart/\n// gramEnd proofErr pair), matching the target OOXML exactly.\n\nconst body = context.document.body;\n\nconst original = \"This is test text with test idea.\";\nconst results = body.search(original, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \" + original);\n}\n\nconst target = results.items[0];\n\n// Flat-OPC WordprocessingML fragment describing the exact replacement\n// run/proofErr layout, inserted via Range.insertOoxml (Office.js's\n// equivalent of Word's Range.InsertXML).\nconst flatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">This is test </w:t></w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r><w:t xml:space=\"preserve\">2  </w:t></w:r>\n            <w:r><w:t>text</w:t></w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> with test idea.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Minor change: insert \"2  \" (a grammar-flagged interjection) right\n# before \"text\" in the second paragraph, turning\n#   \"This is test text with test idea.\"\n# into\n#   \"This is test 2  text with test idea.\"\n# and splitting the run the same way Word's grammar checker would\n# (wrapping the inserted word + the following word in a gramStart/\n# gramEnd proofErr pair), matching the target OOXML exactly.\n\n$d = $word.ActiveDocument\n\n$original = \"This is test text with test idea.\"\n\n$probe = $d.Content\n$found = $probe.Find.Execute($original)\nif (-not $found) {\n    throw \"Target sentence not found: $original\"\n}\n\n# Re-materialise a plain Range over the same bounds - Find's own Range\n# object carries extra search state that makes InsertXML insert instead\n# of replace, so hand InsertXML a fresh Range(start, end) pair.\n$r = $d.Range($probe.Start, $probe.End)\n\n# Flat-OPC WordprocessingML fragment describing the exact replacement\n# run/proofErr layout, applied via Range.InsertXML (the COM twin of\n# Office.js's Range.insertOoxml).\n$flatOpc = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">This is test </w:t></w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r><w:t xml:space=\"preserve\">2  </w:t></w:r>\n            <w:r><w:t>text</w:t></w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> with test idea.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$r.InsertXML($flatOpc)\n"}
